$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A342").Value = "IMX-USD"
$ws.Range("A343").Value = "TAO-USD"
$ws.Range("A344").Value = "GRT-USD"
$ws.Range("A345").Value = "MNT-USD"
